$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Clear the state label (B1, "Texas") and the date stamp (C1, 45531)
# that were merged in by mistake.
$ws.Range("B1:C1").Clear()
